$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = '["https://app-storage-service.pipefy.com/v1/signed/uploads/ebb81f2b-007f-4262-8b9d-5449bed93942/Contrato%20Laura.docx%20-%20Clicksign.pdf?expires_on=1648599380&signature=g4t92KqiRLDSv5pIf3Cq72HFDX5J6BaL%2BaTwgWaDUnQ%3D"]'
$ws.Range("L4").Value = '["https://app-storage-service.pipefy.com/v1/signed/uploads/ecc62a90-29b6-447c-9058-53b38f2a7faf/Contrato%20Gabriel.docx%20-%20Clicksign.pdf?expires_on=1648599380&signature=JLkomHwqOdoY1S8zRFF%2FFi5QWW8ArkMUXmPzFwS2okI%3D"]'
$ws.Range("L5").Value = '["https://app-storage-service.pipefy.com/v1/signed/uploads/d8d41087-022b-43a9-b3a7-c2767d22ee32/Contrato%20Nicolas%20%2B%20Tati%202.docx%20-%20Clicksign.pdf?expires_on=1648599380&signature=9qpZKDD3Kufi0xzg%2FKSqCiSLnpgXyAVqOGpaFCLrLho%3D"]'
$ws.Range("L6").Value = '["https://app-storage-service.pipefy.com/v1/signed/uploads/1d4e7e1b-6181-4d68-9f4c-2af92343a1b0/Contrato%20assinado%20J%C3%A9ssica.pdf?expires_on=1648599380&signature=B21fxsUHieB%2ByWqC6VSBp0eIdaKNZa0RTAz1KUo%2BQg8%3D"]'
$ws.Range("L7").Value = '["https://app-storage-service.pipefy.com/v1/signed/uploads/f034a244-7db2-4690-b6dd-a0e136e82f13/Contrato%20assinado%20Diogo.pdf?expires_on=1648599380&signature=%2F45v5Q%2BREAb3nlLdO7l8RB%2FbYiBtXNjqeaqyD9%2BVOVk%3D"]'
$ws.Range("L8").Value = '["https://app-storage-service.pipefy.com/v1/signed/uploads/f7a08a78-1cd3-4fd4-9e0f-8ffef864bc64/0c6b523f-688c-4083-8149-42b1e8bad674.pdf?expires_on=1648599380&signature=PbU4L1PhCWzZNstKKRzYjqRqEMmkZBEZjM82MQiNLhI%3D"]'
$ws.Range("L9").Value = '["https://app-storage-service.pipefy.com/v1/signed/uploads/ebe95716-517c-4801-90b1-12b470a87eba/d7f66983-f556-4e7a-ab80-c6a9091bd63c.pdf?expires_on=1648599380&signature=iCEmj8S%2BfMJb3guQsqjtIRczlL%2BGaKOiJLwzI3Z45Do%3D"]'
$ws.Range("L10").Value = '["https://app-storage-service.pipefy.com/v1/signed/uploads/4701c047-3d98-4603-9e08-4e07637a1e2b/283e96ef-2df5-4fbd-866a-017ca2a54b46.pdf?expires_on=1648599380&signature=gUjOryNotM%2BDr%2BSYdYokUYwDVxDSUTx7J325eukptJ4%3D"]'
$ws.Range("L11").Value = '["https://app-storage-service.pipefy.com/v1/signed/uploads/b04664f2-a9f4-4126-908b-fcf7d0d927af/0abc913e-411a-432f-a576-8c8ac06efafe.pdf?expires_on=1648599380&signature=zG58mIyHoZxoSichclvRhcUjVMg4tB%2BsjEZZDYM4M1w%3D"]'
$ws.Range("L12").Value = '["https://app-storage-service.pipefy.com/v1/signed/uploads/04631026-fd1b-4df0-922e-b4ae8244b68b/Karinad40db3de-88e0-4ea8-a807-4d8c1b474c99.pdf?expires_on=1648599380&signature=HENxRlHflDPnMtvysRUnAFw7WCczsyJfZkGvW8TifEI%3D"]'
$ws.Range("L13").Value = '["https://app-storage-service.pipefy.com/v1/signed/uploads/0d70b51d-9526-4d01-952e-52987225f2fe/ec5fd910-8712-4b2f-8ec0-5852e10605ae.pdf?expires_on=1648599380&signature=uyH1ZPiBFmc3huEAgVOn1QM6BReuMf8t9GcuiFFOK4E%3D"]'
$ws.Range("L14").Value = '["https://app-storage-service.pipefy.com/v1/signed/uploads/a4645a7b-0477-4878-82bc-8975dd3c019a/0be53e8a-b7e8-41b8-b3e9-db52641599cd.pdf?expires_on=1648599380&signature=kYbhEkRLveKNAmY6zs357be1PdlOdZEy1oXNPi%2Bqnps%3D"]'
$ws.Range("L15").Value = '["https://app-storage-service.pipefy.com/v1/signed/uploads/1cd35f18-267f-4a52-a0a4-a8fcc3a4550b/4243707f-fa7b-480f-9cb8-1b532548f132.pdf?expires_on=1648599380&signature=VsucxZJdbykXaVjRP5rq3FjURmn1s9N%2BKSBtl6E7knI%3D"]'
